# Generate Report for Handback
# Updates the handoff / handback generation timestamps recorded for the
# "ad2609dc-a353-4226-b0f0-a908752174ce.md" file (row 3 of each per-language
# detail table) and the matching "Latest HO Xliff Generate Date" summary
# cell on the Overview sheet.

$wb = $excel.ActiveWorkbook

# --- Overview sheet -------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("G3").Value = "2016-09-06 10:36:29"

# --- zh-cn detail sheet -----------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("H3").Value = "2016-09-06 10:36:16"
$wsZhCn.Range("K3").Value = "2016-09-06 10:37:24"

# --- de-de detail sheet -----------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("H3").Value = "2016-09-06 10:36:29"
$wsDeDe.Range("K3").Value = "2016-09-06 10:37:42"
